$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row
# Cells whose new value would parse as a plain number are first
# formatted as Text so Excel keeps the exact original string
# (matching formatting such as trailing zeros, e.g. "304.46").

$ws.Range("D2").Value = "46.535.52"
$ws.Range("E2").Value = "  +4.87%  "
$ws.Range("D3").Value = "2.295.39"
$ws.Range("E3").Value = "  +2.95%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.46"
$ws.Range("E5").Value = "  +2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.88"
$ws.Range("E6").Value = "  +11.53%  "
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +4.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.44"
$ws.Range("E10").Value = "  +8.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.39"
$ws.Range("E12").Value = "  +5.08%  "
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").Value = "2.645.61"
$ws.Range("E14").Value = "  +3.17%  "
$ws.Range("D15").Value = "2.288.84"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.81"
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.812"
$ws.Range("E17").Value = "  +4.08%  "
$ws.Range("D18").Value = "46.488.24"
$ws.Range("E18").Value = "  +5.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.01"
$ws.Range("E19").Value = "  +7.07%  "
$ws.Range("D20").Value = "0.0₃0940"
$ws.Range("E20").Value = "  +3.40%  "
$ws.Range("E21").Value = "  -0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.91"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "250.26"
$ws.Range("E23").Value = "  +5.86%  "
$ws.Range("E24").Value = "  +2.29%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  +4.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.30"
$ws.Range("E27").Value = "  +7.78%  "
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.88"
$ws.Range("E29").Value = "  +5.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.98"
$ws.Range("E30").Value = "  +3.30%  "
$ws.Range("E31").Value = "  +13.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.61"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "146.71"
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0795"
$ws.Range("E34").Value = "  +3.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.22"
$ws.Range("E35").Value = "  +13.02%  "
$ws.Range("E36").Value = "  +10.54%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.08"
$ws.Range("E38").Value = "  +18.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.77"
$ws.Range("E39").Value = "  +4.65%  "
$ws.Range("E40").Value = "  +9.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.34"
$ws.Range("E41").Value = "  +4.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0302"
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("D45").Value = "1.815.61"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.21"
$ws.Range("E46").Value = "  +19.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.196"
$ws.Range("E47").Value = "  +6.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.11"
$ws.Range("E48").Value = "  +5.58%  "
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "95.22"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "2.517.63"
$ws.Range("E51").Value = "  +2.85%  "
